# Update leve-flip market price data across all job sheets (scheduled runner sync).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 51.57143
$ws.Range("I11").Value = 51.57143
$ws.Range("K11").Value = 51.57143
$ws.Range("M11").Value = 88.42857000000001
$ws.Range("H15").Value = 1806
$ws.Range("I15").Value = 1806
$ws.Range("K15").Value = 5418
$ws.Range("M15").Value = -5249
$ws.Range("H86").Value = 1742.4286
$ws.Range("I86").Value = 1475
$ws.Range("K86").Value = 1475
$ws.Range("M86").Value = -352
$ws.Range("H89").Value = 1742.4286
$ws.Range("I89").Value = 1475
$ws.Range("K89").Value = 7375
$ws.Range("M89").Value = -1759
$ws.Range("I101").Value = 33339128
$ws.Range("J101").Value = 737.5
$ws.Range("K101").Value = 100017384
$ws.Range("L101").Value = 2212.5
$ws.Range("M101").Value = -100015762
$ws.Range("N101").Value = -5456.5
$ws.Range("H135").Value = 4474.25
$ws.Range("I135").Value = 3605
$ws.Range("K135").Value = 32445
$ws.Range("M135").Value = -29910

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 100.5
$ws.Range("I5").Value = 100.5
$ws.Range("K5").Value = 100.5
$ws.Range("M5").Value = 11.5
$ws.Range("H32").Value = 4281729.5
$ws.Range("I32").Value = 4121831.5
$ws.Range("K32").Value = 4121831.5
$ws.Range("M32").Value = -4121544.5
$ws.Range("H50").Value = 4447.375
$ws.Range("I50").Value = 6783.2
$ws.Range("K50").Value = 6783.2
$ws.Range("M50").Value = -6069.2
$ws.Range("H101").Value = 50301
$ws.Range("J101").Value = 50301
$ws.Range("L101").Value = 50301
$ws.Range("N101").Value = -56791
$ws.Range("H113").Value = 89990
$ws.Range("J113").Value = 89990
$ws.Range("L113").Value = 89990
$ws.Range("N113").Value = -98668

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100.5
$ws.Range("I4").Value = 100.5
$ws.Range("K4").Value = 100.5
$ws.Range("M4").Value = 14.5
$ws.Range("H22").Value = 165.8
$ws.Range("I22").Value = 150.71428
$ws.Range("J22").Value = 201
$ws.Range("K22").Value = 150.71428
$ws.Range("L22").Value = 201
$ws.Range("M22").Value = 22.28572
$ws.Range("N22").Value = -547
$ws.Range("H80").Value = 253.5
$ws.Range("I80").Value = 260
$ws.Range("J80").Value = 247
$ws.Range("K80").Value = 260
$ws.Range("L80").Value = 247
$ws.Range("M80").Value = 738
$ws.Range("N80").Value = -2243
$ws.Range("H83").Value = 253.5
$ws.Range("I83").Value = 260
$ws.Range("J83").Value = 247
$ws.Range("K83").Value = 1300
$ws.Range("L83").Value = 1235
$ws.Range("M83").Value = 3692
$ws.Range("N83").Value = -11219
$ws.Range("H86").Value = 1596.1578
$ws.Range("I86").Value = 1573.25
$ws.Range("J86").Value = 1635.4286
$ws.Range("K86").Value = 1573.25
$ws.Range("L86").Value = 1635.4286
$ws.Range("M86").Value = -450.25
$ws.Range("N86").Value = -3881.4286
$ws.Range("H89").Value = 1596.1578
$ws.Range("I89").Value = 1573.25
$ws.Range("J89").Value = 1635.4286
$ws.Range("K89").Value = 7866.25
$ws.Range("L89").Value = 8177.143
$ws.Range("M89").Value = -2250.25
$ws.Range("N89").Value = -19409.143

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2533.3333
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 2533.3333
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H99").Value = 1711.9166
$ws.Range("I99").Value = 1532.3334
$ws.Range("J99").Value = 2250.6667
$ws.Range("K99").Value = 1532.3334
$ws.Range("L99").Value = 2250.6667
$ws.Range("M99").Value = -34.33339999999998
$ws.Range("N99").Value = -5246.6667
$ws.Range("H126").Value = 1711.9166
$ws.Range("I126").Value = 1532.3334
$ws.Range("J126").Value = 2250.6667
$ws.Range("K126").Value = 4597.0002
$ws.Range("L126").Value = 6752.000100000001
$ws.Range("M126").Value = -2127.0002
$ws.Range("N126").Value = -11692.0001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2380
$ws.Range("J81").Value = 2380
$ws.Range("L81").Value = 7140
$ws.Range("N81").Value = -9386
$ws.Range("H84").Value = 2380
$ws.Range("J84").Value = 2380
$ws.Range("L84").Value = 21420
$ws.Range("N84").Value = -32652
$ws.Range("H109").Value = 1956.75
$ws.Range("I109").Value = 1675.6666
$ws.Range("K109").Value = 5026.9998
$ws.Range("M109").Value = -3986.9998
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H112").Value = 16277.375
$ws.Range("J112").Value = 18285.285
$ws.Range("L112").Value = 54855.855
$ws.Range("N112").Value = -57071.855
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("H119").Value = 1495
$ws.Range("I119").Value = 1495
$ws.Range("K119").Value = 4485
$ws.Range("M119").Value = 353
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("H122").Value = 600
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2108.1667
$ws.Range("J132").Value = 1225
$ws.Range("L132").Value = 11025
$ws.Range("N132").Value = -16085

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3165.6667
$ws.Range("I80").Value = 3165.6667
$ws.Range("K80").Value = 3165.6667
$ws.Range("M80").Value = -2167.6667
$ws.Range("H83").Value = 3165.6667
$ws.Range("I83").Value = 3165.6667
$ws.Range("K83").Value = 15828.3335
$ws.Range("M83").Value = -10836.3335
$ws.Range("H122").Value = 5279.8
$ws.Range("I122").Value = 6399
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 19197
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -16747
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 6020.926
$ws.Range("I132").Value = 6273.5415
$ws.Range("K132").Value = 18820.6245
$ws.Range("M132").Value = -16290.6245

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 14999
$ws.Range("J20").Value = 14999
$ws.Range("L20").Value = 14999
$ws.Range("N20").Value = -15451
$ws.Range("H22").Value = 1945.5625
$ws.Range("I22").Value = 2030.6428
$ws.Range("J22").Value = 1350
$ws.Range("K22").Value = 2030.6428
$ws.Range("L22").Value = 1350
$ws.Range("M22").Value = -1735.6428
$ws.Range("N22").Value = -1940
$ws.Range("H27").Value = 1945.5625
$ws.Range("I27").Value = 2030.6428
$ws.Range("J27").Value = 1350
$ws.Range("K27").Value = 2030.6428
$ws.Range("L27").Value = 1350
$ws.Range("M27").Value = -1923.6428
$ws.Range("N27").Value = -1564
$ws.Range("H40").Value = 4350
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5272

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 32125.5
$ws.Range("I62").Value = 24500.666
$ws.Range("K62").Value = 24500.666
$ws.Range("M62").Value = -23876.666
$ws.Range("H65").Value = 32125.5
$ws.Range("I65").Value = 24500.666
$ws.Range("K65").Value = 122503.33
$ws.Range("M65").Value = -119383.33
$ws.Range("H81").Value = 3334450
$ws.Range("I81").Value = 1674
$ws.Range("J81").Value = 10000002
$ws.Range("K81").Value = 3348
$ws.Range("L81").Value = 20000004
$ws.Range("M81").Value = -2287
$ws.Range("N81").Value = -20002126
$ws.Range("H84").Value = 3334450
$ws.Range("I84").Value = 1674
$ws.Range("J84").Value = 10000002
$ws.Range("K84").Value = 16740
$ws.Range("L84").Value = 100000020
$ws.Range("M84").Value = -11436
$ws.Range("N84").Value = -100010628
$ws.Range("H113").Value = 520.5714
$ws.Range("J113").Value = 581
$ws.Range("L113").Value = 1743
$ws.Range("N113").Value = -6083
$ws.Range("H122").Value = 9812
$ws.Range("J122").Value = 9828.833000000001
$ws.Range("L122").Value = 29486.499
$ws.Range("N122").Value = -34386.499
$ws.Range("H136").Value = 937.4286
$ws.Range("J136").Value = 949.5714
$ws.Range("L136").Value = 2848.7142
$ws.Range("N136").Value = -7948.7142
